$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$s.MoveTo(19)
